# Insert a new weekly price record for "Femacal de La Calera" / Zapallo / Camote
# just above the current row 303 (everything from row 303 down shifts to +1),
# matching the "Fruta / hortaliza, semanal" refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 303..351 down to 304..352, inheriting formatting from
# the row above (same behaviour as Excel's own "Insert Sheet Rows").
$ws.Rows("303:303").Insert()

# Populate the newly opened row 303 with the new record.
$ws.Range("A303").Value = 3
$ws.Range("B303").Value = "Femacal de La Calera"
$ws.Range("C303").Value = "Coquimbo"
$ws.Range("D303").Value = 44505
$ws.Range("E303").Value = 5
$ws.Range("F303").Value = 100112045
$ws.Range("G303").Value = "Zapallo"
$ws.Range("H303").Value = "Camote"
$ws.Range("I303").Value = "1a (guarda)"
$ws.Range("J303").Value = 160
$ws.Range("K303").Value = 550
$ws.Range("L303").Value = 550
$ws.Range("M303").Value = 550
$ws.Range("N303").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O303").Value = "Provincia de Talca"
$ws.Range("P303").Value = 550
$ws.Range("Q303").Value = 1
$ws.Range("R303").Value = "Hortaliza"
